# Regenerate save_data to use K (strikeouts) values instead of Strike# values
# in column G. Other columns remain unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 5
    5  = 2
    6  = 4
    7  = 4
    8  = 6
    9  = 2
    10 = 1
    11 = 2
    12 = 4
    13 = 2
    14 = 3
    15 = 2
    16 = 3
    17 = 2
    18 = 6
    19 = 3
    20 = 2
    21 = 2
    22 = 2
    23 = 0
    24 = 4
    25 = 4
    26 = 4
    27 = 2
    28 = 4
    29 = 4
    30 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
